$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Industries") values for rows 26-176 change from 1 to 0
$ws.Range("H26:H176").Value = 0
